# kibon-3327: improve statistik vorlage - avoid reference error in excel formula
#
# The template row (row 11) contains formulas that reference placeholder
# cells (J11, Y11, AK11, ...) which are not yet real numbers/dates. That
# causes EOMONTH()/multiplication to blow up with #VALUE! errors in the
# raw template. Wrap the formulas with ISNUMBER() guards so they
# gracefully evaluate to "" instead of erroring out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M11").Formula = '=IF(ISNUMBER(J11),IF(Y11<=EOMONTH(J11,12),"X",""),"")'
$ws.Range("Q11").Formula = '=IF(ISNUMBER(J11),IF(Y11<=EOMONTH(J11,12),"X",""),"")'
$ws.Range("R11").Formula = '=IF(ISNUMBER(J11),IF(AND(Y11>=EOMONTH(J11,13),Y11<=EOMONTH(J11,48)),"X",""),"")'
$ws.Range("S11").Formula = '=IF(ISNUMBER(J11),IF(AND(Y11>=EOMONTH(J11,48),Y11<=EOMONTH(J11,72)),"X",""),"")'
$ws.Range("T11").Formula = '=IF(ISNUMBER(J11),IF(Y11>=EOMONTH(J11,73),"X",""),"")'
$ws.Range("AN11").Formula = '=IF(ISNUMBER(AK11),AK11*AD11,"")'

# Reflect the author's scroll/selection state: the view was scrolled so
# column D is the left-most visible column, with the active cell on J11.
$ws.Application.Goto($ws.Range("J11"), $false)
$excel.ActiveWindow.ScrollColumn = $ws.Range("D1").Column
$ws.Range("J11").Select()
